# Auto update Excel log
# Appends newly captured sensor-log rows to the ALERTS and mmWave sheets.

$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: append 4 new FALL_DETECTED rows (10-13) ---
$alerts = $wb.Worksheets.Item("ALERTS")

$alertRows = @(
    @{ Row=10; Date="2026-02-01"; Timestamp="00:24:31"; Hour="00:00"; Location="Living Room"; Value="CRITICAL"; Status="FALL_DETECTED" },
    @{ Row=11; Date="2026-02-01"; Timestamp="00:24:36"; Hour="00:00"; Location="Living Room"; Value="CRITICAL"; Status="FALL_DETECTED" },
    @{ Row=12; Date="2026-02-01"; Timestamp="00:24:39"; Hour="00:00"; Location="Living Room"; Value="CRITICAL"; Status="FALL_DETECTED" },
    @{ Row=13; Date="2026-02-01"; Timestamp="00:24:43"; Hour="00:00"; Location="Living Room"; Value="CRITICAL"; Status="FALL_DETECTED" }
)

foreach ($r in $alertRows) {
    # Column A holds ISO-style dates (e.g. 2026-02-01). Force Text format
    # first so Excel doesn't auto-convert the literal into a date serial
    # number, keeping it a plain string like the rest of the log.
    $alerts.Cells.Item($r.Row, 1).NumberFormat = "@"
    $alerts.Cells.Item($r.Row, 1).Value = $r.Date
    $alerts.Cells.Item($r.Row, 2).Value = $r.Timestamp
    $alerts.Cells.Item($r.Row, 3).Value = $r.Hour
    $alerts.Cells.Item($r.Row, 4).Value = $r.Location
    $alerts.Cells.Item($r.Row, 5).Value = $r.Value
    $alerts.Cells.Item($r.Row, 6).Value = $r.Status
}

# --- mmWave sheet: append 3 new motion rows (56-58) ---
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwaveRows = @(
    @{ Row=56; Date="2026-02-01"; Timestamp="00:24:55"; Hour="00:00"; Location="Living Room"; Value="NO_MOTION_DETECTED"; Status="Inactive" },
    @{ Row=57; Date="2026-02-01"; Timestamp="00:25:05"; Hour="00:00"; Location="Living Room"; Value="PRESENCE_DETECTED"; Status="Active" },
    @{ Row=58; Date="2026-02-01"; Timestamp="00:25:16"; Hour="00:00"; Location="Living Room"; Value="PRESENCE_DETECTED"; Status="Active" }
)

foreach ($r in $mmwaveRows) {
    $mmwave.Cells.Item($r.Row, 1).NumberFormat = "@"
    $mmwave.Cells.Item($r.Row, 1).Value = $r.Date
    $mmwave.Cells.Item($r.Row, 2).Value = $r.Timestamp
    $mmwave.Cells.Item($r.Row, 3).Value = $r.Hour
    $mmwave.Cells.Item($r.Row, 4).Value = $r.Location
    $mmwave.Cells.Item($r.Row, 5).Value = $r.Value
    $mmwave.Cells.Item($r.Row, 6).Value = $r.Status
}
